$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 280 — this shifts the existing rows
# 280..306 down to 281..307 (and carries the row-280 formatting, e.g. the
# date number format on column D, down with them as Excel normally does).
$ws.Rows("280:280").Insert()

# Populate the newly inserted row 280 with the new weekly price record.
$ws.Range("A280").Value = 3
$ws.Range("B280").Value = "Femacal de La Calera"
$ws.Range("C280").Value = "Coquimbo"
$ws.Range("D280").Value = 45194
$ws.Range("E280").Value = 5
$ws.Range("F280").Value = 100112026
$ws.Range("G280").Value = "Haba"
$ws.Range("H280").Value = "Sin especificar"
$ws.Range("I280").Value = "Primera"
$ws.Range("J280").Value = 75
$ws.Range("K280").Value = 13000
$ws.Range("L280").Value = 14000
$ws.Range("M280").Value = 13467
$ws.Range("N280").Value = "$/saco 25 kilos"
$ws.Range("O280").Value = "Provincia de Limarí"
$ws.Range("P280").Value = 539
$ws.Range("Q280").Value = 25
$ws.Range("R280").Value = "Hortaliza"
